$wb = $excel.ActiveWorkbook

# --- OpsTracker sheet updates ---
$ops = $wb.Worksheets.Item("OpsTracker")

# Item 15 "Stamp Pad": Status Done -> WIP
$ops.Range("D16").Value = "WIP"

# Item 27 "Merch for Marketing Team": add a Comment
$ops.Range("E28").Value = "Return of Merch (Debashish and Sayan)"

# --- InternalAdmin sheet updates ---
$admin = $wb.Worksheets.Item("InternalAdmin")

# Furniture Repairing: Todo -> WIP
$admin.Range("C2").Value = "WIP"

# AC Service: Todo -> WIP
$admin.Range("C3").Value = "WIP"

# Banking issue: Todo -> WIP
$admin.Range("C6").Value = "WIP"

# Pujo Subscription: Todo -> Done
$admin.Range("C9").Value = "Done"

# ID card jacket purchase: Owner Victor -> Anirban
$admin.Range("D10").Value = "Anirban"

# New SIM: Todo -> Done
$admin.Range("C12").Value = "Done"
